$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1.141346
    "H2" = 3.424038
    "I2" = 0.07840323688201205
    "J2" = 0.07840323688201206
    "M2" = 6.673575666666667
    "N2" = 20.020727
    "O2" = 0.1592723389521969
    "P2" = 0.1592723389521969
    "Q2" = 7.616858892847334
    "R2" = 68.55173003562601
    "S2" = 0.01248746691962121
    "T2" = 0.01248746691962121
    "G3" = 1.141346
    "H3" = 3.424038
    "I3" = 0.07840323688201205
    "J3" = 0.07840323688201206
    "O3" = 0.683825744632005
    "P3" = 0.6838257446320051
    "Q3" = 32.70250338774466
    "R3" = 294.322530489702
    "S3" = 0.05361415184240136
    "T3" = 0.05361415184240138
    "G4" = 1.141346
    "H4" = 3.424038
    "I4" = 0.07840323688201205
    "J4" = 0.07840323688201206
    "K4" = 2
    "L4" = 0.6666666666666666
    "M4" = 0.2748103333333333
    "N4" = 0.824431
    "O4" = 0.006558655620982129
    "P4" = 0.006558655620982129
    "Q4" = 0.3136536747086667
    "R4" = 2.822883072378
    "S4" = 0.0005142198302794017
    "T4" = 0.0005142198302794018
    "G5" = 1.141346
    "H5" = 3.424038
    "I5" = 0.07840323688201205
    "J5" = 0.07840323688201206
    "M5" = 5.865251333333333
    "N5" = 17.595754
    "O5" = 0.1399807756834942
    "P5" = 0.1399807756834942
    "Q5" = 6.694281148294666
    "R5" = 60.248530334652
    "S5" = 0.01097494591484078
    "T5" = 0.01097494591484079
    "G6" = 1.141346
    "H6" = 3.424038
    "I6" = 0.07840323688201205
    "J6" = 0.07840323688201206
    "M6" = 0.4341923333333333
    "N6" = 1.302577
    "O6" = 0.01036248511132167
    "P6" = 0.01036248511132167
    "Q6" = 0.4955636828806666
    "R6" = 4.460073145926
    "S6" = 0.000812452374869276
    "T6" = 0.0008124523748692763
    "I7" = 0.118601898037951
    "J7" = 0.118601898037951
    "M7" = 6.673575666666667
    "N7" = 20.020727
    "O7" = 0.1592723389521969
    "P7" = 0.1592723389521969
    "Q7" = 11.52215084102222
    "R7" = 103.6993575692
    "S7" = 0.01889000170467443
    "T7" = 0.01889000170467444
    "I8" = 0.118601898037951
    "J8" = 0.118601898037951
    "O8" = 0.683825744632005
    "P8" = 0.6838257446320051
    "S8" = 0.08110303124057097
    "T8" = 0.081103031240571
    "I9" = 0.118601898037951
    "J9" = 0.118601898037951
    "K9" = 2
    "L9" = 0.6666666666666666
    "M9" = 0.2748103333333333
    "N9" = 0.824431
    "O9" = 0.006558655620982129
    "P9" = 0.006558655620982129
    "Q9" = 0.4744692008444444
    "R9" = 4.270222807600001
    "S9" = 0.0007778690052257566
    "T9" = 0.0007778690052257568
    "I10" = 0.118601898037951
    "J10" = 0.118601898037951
    "M10" = 5.865251333333333
    "N10" = 17.595754
    "O10" = 0.1399807756834942
    "P10" = 0.1399807756834942
    "Q10" = 10.12655193537778
    "R10" = 91.13896741840001
    "S10" = 0.01660198568488706
    "T10" = 0.01660198568488707
    "I11" = 0.118601898037951
    "J11" = 0.118601898037951
    "M11" = 0.4341923333333333
    "N11" = 1.302577
    "O11" = 0.01036248511132167
    "P11" = 0.01036248511132167
    "Q11" = 0.7496475365777777
    "R11" = 6.7468278292
    "S11" = 0.001229010402592758
    "T11" = 0.001229010402592759
    "G12" = 4.821393333333334
    "H12" = 14.46418
    "I12" = 0.331199166260439
    "J12" = 0.3311991662604391
    "M12" = 6.673575666666667
    "N12" = 20.020727
    "O12" = 0.1592723389521969
    "P12" = 0.1592723389521969
    "Q12" = 32.17593322876223
    "R12" = 289.58339905886
    "S12" = 0.05275086586931768
    "T12" = 0.05275086586931769
    "G13" = 4.821393333333334
    "H13" = 14.46418
    "I13" = 0.331199166260439
    "J13" = 0.3311991662604391
    "O13" = 0.683825744632005
    "P13" = 0.6838257446320051
    "Q13" = 138.1453405163578
    "R13" = 1243.30806464722
    "S13" = 0.226482516489544
    "T13" = 0.226482516489544
    "G14" = 4.821393333333334
    "H14" = 14.46418
    "I14" = 0.331199166260439
    "J14" = 0.3311991662604391
    "K14" = 2
    "L14" = 0.6666666666666666
    "M14" = 0.2748103333333333
    "N14" = 0.824431
    "O14" = 0.006558655620982129
    "P14" = 0.006558655620982129
    "Q14" = 1.324968709064445
    "R14" = 11.92471838158
    "S14" = 0.002172221273458623
    "T14" = 0.002172221273458624
    "G15" = 4.821393333333334
    "H15" = 14.46418
    "I15" = 0.331199166260439
    "J15" = 0.3311991662604391
    "M15" = 5.865251333333333
    "N15" = 17.595754
    "O15" = 0.1399807756834942
    "P15" = 0.1399807756834942
    "Q15" = 28.27868367685778
    "R15" = 254.50815309172
    "S15" = 0.0463615161988628
    "T15" = 0.04636151619886282
    "G16" = 4.821393333333334
    "H16" = 14.46418
    "I16" = 0.331199166260439
    "J16" = 0.3311991662604391
    "M16" = 0.4341923333333333
    "N16" = 1.302577
    "O16" = 0.01036248511132167
    "P16" = 0.01036248511132167
    "Q16" = 2.093412021317778
    "R16" = 18.84070819186
    "S16" = 0.00343204642925595
    "T16" = 0.003432046429255952
    "G17" = 1.935744333333333
    "H17" = 5.807233
    "I17" = 0.1329733678563256
    "J17" = 0.1329733678563257
    "M17" = 6.673575666666667
    "N17" = 20.020727
    "O17" = 0.1592723389521969
    "P17" = 0.1592723389521969
    "Q17" = 12.91833627982122
    "R17" = 116.265026518391
    "S17" = 0.02117897931682787
    "T17" = 0.02117897931682787
    "G18" = 1.935744333333333
    "H18" = 5.807233
    "I18" = 0.1329733678563256
    "J18" = 0.1329733678563257
    "O18" = 0.683825744632005
    "P18" = 0.6838257446320051
    "Q18" = 55.46406227265077
    "R18" = 499.176560453857
    "S18" = 0.09093061229057739
    "T18" = 0.09093061229057742
    "G19" = 1.935744333333333
    "H19" = 5.807233
    "I19" = 0.1329733678563256
    "J19" = 0.1329733678563257
    "K19" = 2
    "L19" = 0.6666666666666666
    "M19" = 0.2748103333333333
    "N19" = 0.824431
    "O19" = 0.006558655620982129
    "P19" = 0.006558655620982129
    "Q19" = 0.5319625454914444
    "R19" = 4.787662909423
    "S19" = 0.0008721265265318144
    "T19" = 0.0008721265265318147
    "G20" = 1.935744333333333
    "H20" = 5.807233
    "I20" = 0.1329733678563256
    "J20" = 0.1329733678563257
    "M20" = 5.865251333333333
    "N20" = 17.595754
    "O20" = 0.1399807756834942
    "P20" = 0.1399807756834942
    "Q20" = 11.35362703207578
    "R20" = 102.182643288682
    "S20" = 0.01861371517777507
    "T20" = 0.01861371517777508
    "G21" = 1.935744333333333
    "H21" = 5.807233
    "I21" = 0.1329733678563256
    "J21" = 0.1329733678563257
    "M21" = 0.4341923333333333
    "N21" = 1.302577
    "O21" = 0.01036248511132167
    "P21" = 0.01036248511132167
    "Q21" = 0.8404853488267777
    "R21" = 7.564368139440999
    "S21" = 0.001377934544613474
    "T21" = 0.001377934544613475
    "G22" = 4.932366666666667
    "H22" = 14.7971
    "I22" = 0.3388223309632722
    "J22" = 0.3388223309632722
    "M22" = 6.673575666666667
    "N22" = 20.020727
    "O22" = 0.1592723389521969
    "P22" = 0.1592723389521969
    "Q22" = 32.91652216574445
    "R22" = 296.2486994917
    "S22" = 0.05396502514175575
    "T22" = 0.05396502514175575
    "G23" = 4.932366666666667
    "H23" = 14.7971
    "I23" = 0.3388223309632722
    "J23" = 0.3388223309632722
    "O23" = 0.683825744632005
    "P23" = 0.6838257446320051
    "Q23" = 138.1453405163578
    "R23" = 1271.9251117859
    "S23" = 0.2316954327689112
    "T23" = 0.2316954327689113
    "G24" = 4.932366666666667
    "H24" = 14.7971
    "I24" = 0.3388223309632722
    "J24" = 0.3388223309632722
    "K24" = 2
    "L24" = 0.6666666666666666
    "M24" = 0.2748103333333333
    "N24" = 0.824431
    "O24" = 0.006558655620982129
    "P24" = 0.006558655620982129
    "Q24" = 1.355465327788889
    "R24" = 12.1991879501
    "S24" = 0.002222218985486532
    "T24" = 0.002222218985486533
    "G25" = 4.932366666666667
    "H25" = 14.7971
    "I25" = 0.3388223309632722
    "J25" = 0.3388223309632722
    "M25" = 5.865251333333333
    "N25" = 17.595754
    "O25" = 0.1399807756834942
    "P25" = 0.1399807756834942
    "Q25" = 28.92957016815556
    "R25" = 260.3661315134
    "S25" = 0.04742861270712842
    "T25" = 0.04742861270712843
    "G26" = 4.932366666666667
    "H26" = 14.7971
    "I26" = 0.3388223309632722
    "J26" = 0.3388223309632722
    "M26" = 0.4341923333333333
    "N26" = 1.302577
    "O26" = 0.01036248511132167
    "P26" = 0.01036248511132167
    "Q26" = 2.141595791855555
    "R26" = 19.2743621267
    "S26" = 0.003511041359990212
    "T26" = 0.003511041359990213
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
